$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the pinmapping mistake: swap PTC12/PTC4 between D4 and D5 ---
$ws.Range("D4").Value = "PTC4"
$ws.Range("D5").Value = "PTC12"

# --- Remove the stray leftover value in I11 (row becomes empty and disappears) ---
$ws.Range("I11").ClearContents()

# --- Enter the new "lcd" pin-mapping block content values, in the exact
#     order the original author typed them (this controls the order new
#     entries are appended to the shared-strings table) ---
$ws.Range("A19").Value = "lcd"
$ws.Range("B20").Value = "nCS"
$ws.Range("C20").Value = "D10"
$ws.Range("C21").Value = "D7"
$ws.Range("B21").Value = "A0(miso?)"
$ws.Range("B22").Value = "SCK"
$ws.Range("C22").Value = "D13"
$ws.Range("B23").Value = "reset"
$ws.Range("C23").Value = "D12"
$ws.Range("B24").Value = "mosi"
$ws.Range("C24").Value = "D11"
$ws.Range("D20").Value = "PTD0"
$ws.Range("D21").Value = "PTC3"
$ws.Range("D22").Value = "PTD1"
$ws.Range("D23").Value = "PTD3"
$ws.Range("D24").Value = "PTD2"
$ws.Range("C19").Value = "AP"
$ws.Range("D19").Value = "FRDM"

# --- Enter the new "speaker" pin-mapping block content values ---
$ws.Range("A26").Value = "speaker"
$ws.Range("B27").Value = "speaker"
$ws.Range("C27").Value = "D6"
$ws.Range("D27").Value = "PTC2"
$ws.Range("C26").Value = "AP"
$ws.Range("D26").Value = "FRDM"

# --- Apply correct cell formatting to the new cells by copying formats
#     from existing cells that already use the desired style ---
$ws.Range("A12").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A12").Copy()
$ws.Range("A26").PasteSpecial(-4122)

$ws.Range("C3").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("B26").PasteSpecial(-4122)

$ws.Range("C12").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("D12").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("D12").Copy()
$ws.Range("D26").PasteSpecial(-4122)

$ws.Range("B8").Copy()
$ws.Range("B20:B24").PasteSpecial(-4122)
$ws.Range("B8").Copy()
$ws.Range("B27").PasteSpecial(-4122)

$ws.Range("C3").Copy()
$ws.Range("C20:C24").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("D20:D24").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("D27").PasteSpecial(-4122)

# B19/B26 are blank placeholder cells (formatting only, no value)
$ws.Range("B19").ClearContents()
$ws.Range("B26").ClearContents()

# --- Update active cell selection to match target ---
$ws.Range("D5").Select()
